$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.828.73'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '2.281.64'
$ws.Range("E3").Value = '  +1.42%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = "'251.04"
$ws.Range("E5").Value = '  +0.71%  '
$ws.Range("D6").Value = "'0.644"
$ws.Range("E6").Value = '  +1.83%  '
$ws.Range("D7").Value = "'75.17"
$ws.Range("E7").Value = '  +6.58%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = "'0.640"
$ws.Range("E9").Value = '  -2.93%  '
$ws.Range("D10").Value = "'39.66"
$ws.Range("E10").Value = '  +2.72%  '
$ws.Range("D11").Value = "'0.0978"
$ws.Range("E11").Value = '  +1.75%  '
$ws.Range("D12").Value = "'7.43"
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("E13").Value = '  +1.99%  '
$ws.Range("D14").Value = '2.625.51'
$ws.Range("E14").Value = '  +1.79%  '
$ws.Range("D15").Value = "'15.07"
$ws.Range("E15").Value = '  +2.84%  '
$ws.Range("D16").Value = "'0.868"
$ws.Range("E16").Value = '  -1.10%  '
$ws.Range("D17").Value = '2.279.78'
$ws.Range("E17").Value = '  +1.57%  '
$ws.Range("D18").Value = '42.716.89'
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("D19").Value = '0.0₃0999'
$ws.Range("E19").Value = '  +1.36%  '
$ws.Range("D20").Value = "'6.22"
$ws.Range("E20").Value = '  -0.60%  '
$ws.Range("D21").Value = "'72.45"
$ws.Range("E21").Value = '  -0.53%  '
$ws.Range("D22").Value = "'236.60"
$ws.Range("E22").Value = '  +1.41%  '
$ws.Range("E23").Value = '  +5.55%  '
$ws.Range("E24").Value = '  -1.70%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = "'11.30"
$ws.Range("E26").Value = '  -0.96%  '
$ws.Range("D27").Value = "'2.40"
$ws.Range("E27").Value = '  -0.89%  '
$ws.Range("E28").Value = '  +1.64%  '
$ws.Range("D29").Value = "'167.46"
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("D30").Value = "'21.06"
$ws.Range("E30").Value = '  +1.22%  '
$ws.Range("D31").Value = "'0.0870"
$ws.Range("E31").Value = '  +9.32%  '
$ws.Range("D32").Value = "'6.39"
$ws.Range("E32").Value = '  -1.13%  '
$ws.Range("D33").Value = "'0.125"
$ws.Range("E33").Value = '  +0.26%  '
$ws.Range("D34").Value = "'32.00"
$ws.Range("E34").Value = '  +2.34%  '
$ws.Range("E35").Value = '  +2.45%  '
$ws.Range("D36").Value = "'4.56"
$ws.Range("E36").Value = '  +3.48%  '
$ws.Range("D37").Value = "'4.77"
$ws.Range("E37").Value = '  +1.61%  '
$ws.Range("D38").Value = "'0.0306"
$ws.Range("E38").Value = '  -5.15%  '
$ws.Range("D39").Value = "'13.58"
$ws.Range("E39").Value = '  +10.19%  '
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("D41").Value = "'5.89"
$ws.Range("E41").Value = '  +1.96%  '
$ws.Range("E42").Value = '  +3.79%  '
$ws.Range("D43").Value = "'61.82"
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("D44").Value = "'8.95"
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("D45").Value = "'106.38"
$ws.Range("E45").Value = '  +13.11%  '
$ws.Range("D46").Value = "'4.73"
$ws.Range("E46").Value = '  -3.16%  '
$ws.Range("E47").Value = '  -0.57%  '
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("E49").Value = '  +0.23%  '
$ws.Range("E50").Value = '  -1.29%  '
$ws.Range("D51").Value = "'4.22"
$ws.Range("E51").Value = '  +0.30%  '
